# Backup QR Scanner data - rename sheet + append new scan/manual log rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Session" to "Neurology"
$ws.Name = "Neurology"

# New log rows to append (rows 37-39), matching the existing text-stored layout
$newRows = @(
    @{ Row = 37; A = "190333"; B = "Neurology"; C = "16/12/2025"; D = "10:13:46"; E = "Manual"; F = "emp17.farah.a.youssef@gmail.com" },
    @{ Row = 38; A = "191007"; B = "Neurology"; C = "16/12/2025"; D = "10:16:24"; E = "Scan";   F = "emp17.farah.a.youssef@gmail.com" },
    @{ Row = 39; A = "202051"; B = "Neurology"; C = "16/12/2025"; D = "10:48:43"; E = "Scan";   F = "emp17.farah.a.youssef@gmail.com" }
)

# Every column in this log is kept as text (the sheet already flags
# numberStoredAsText for A1:F36), so force text format before writing values
# to avoid the numeric-looking IDs being stored as numbers.
$ws.Range("A37:F39").NumberFormat = "@"

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
}
